$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Update the "Recorded By" (column G) values: reorder the comma-separated
# list of recorder names/emails for the specific rows that changed.
$ws.Range("G2").Value = 'backup@backdoor.com, System, system'
$ws.Range("G3").Value = 'System, dnasr281@gmail.com'
$ws.Range("G4").Value = 'backup@backdoor.com, System'
$ws.Range("G5").Value = 'backup@backdoor.com, System'
$ws.Range("G6").Value = 'System, dnasr281@gmail.com'
$ws.Range("G7").Value = 'System, admin@admin.com'
$ws.Range("G8").Value = 'backup@backdoor.com, System'
$ws.Range("G10").Value = 'System, dnasr281@gmail.com'
$ws.Range("G12").Value = 'System, dnasr281@gmail.com'
$ws.Range("G13").Value = 'System, dnasr281@gmail.com'
$ws.Range("G14").Value = 'System, dnasr281@gmail.com'
$ws.Range("G15").Value = 'System, dnasr281@gmail.com'
$ws.Range("G18").Value = 'System, dnasr281@gmail.com'
$ws.Range("G19").Value = 'System, dnasr281@gmail.com'
$ws.Range("G20").Value = 'System, dnasr281@gmail.com'
$ws.Range("G21").Value = 'System, dnasr281@gmail.com'
$ws.Range("G22").Value = 'System, dnasr281@gmail.com'
$ws.Range("G29").Value = 'backup@backdoor.com, System, system'
$ws.Range("G30").Value = 'System, dnasr281@gmail.com'
$ws.Range("G31").Value = 'backup@backdoor.com, System'
$ws.Range("G32").Value = 'backup@backdoor.com, System'
$ws.Range("G33").Value = 'System, dnasr281@gmail.com'
$ws.Range("G34").Value = 'System, admin@admin.com'
$ws.Range("G35").Value = 'backup@backdoor.com, System'
$ws.Range("G37").Value = 'System, dnasr281@gmail.com'
$ws.Range("G39").Value = 'System, dnasr281@gmail.com'
$ws.Range("G40").Value = 'System, dnasr281@gmail.com'
$ws.Range("G41").Value = 'System, dnasr281@gmail.com'
$ws.Range("G42").Value = 'System, dnasr281@gmail.com'
$ws.Range("G45").Value = 'System, dnasr281@gmail.com'
$ws.Range("G46").Value = 'System, dnasr281@gmail.com'
$ws.Range("G47").Value = 'System, dnasr281@gmail.com'
$ws.Range("G48").Value = 'System, dnasr281@gmail.com'
$ws.Range("G49").Value = 'System, dnasr281@gmail.com'
$ws.Range("G56").Value = 'backup@backdoor.com, System, system'
$ws.Range("G57").Value = 'System, dnasr281@gmail.com'
$ws.Range("G58").Value = 'backup@backdoor.com, System'
$ws.Range("G59").Value = 'backup@backdoor.com, System'
$ws.Range("G60").Value = 'System, dnasr281@gmail.com'
$ws.Range("G61").Value = 'System, admin@admin.com'
$ws.Range("G62").Value = 'backup@backdoor.com, System'
$ws.Range("G64").Value = 'System, dnasr281@gmail.com'
$ws.Range("G66").Value = 'System, dnasr281@gmail.com'
$ws.Range("G67").Value = 'System, dnasr281@gmail.com'
$ws.Range("G68").Value = 'System, dnasr281@gmail.com'
$ws.Range("G69").Value = 'System, dnasr281@gmail.com'
$ws.Range("G72").Value = 'System, dnasr281@gmail.com'
$ws.Range("G73").Value = 'System, dnasr281@gmail.com'
$ws.Range("G74").Value = 'System, dnasr281@gmail.com'
$ws.Range("G75").Value = 'System, dnasr281@gmail.com'
$ws.Range("G76").Value = 'System, dnasr281@gmail.com'
$ws.Range("G83").Value = 'backup@backdoor.com, System'
$ws.Range("G84").Value = 'backup@backdoor.com, System'
$ws.Range("G85").Value = 'backup@backdoor.com, System'
$ws.Range("G86").Value = 'System, dnasr281@gmail.com'
$ws.Range("G87").Value = 'System, dnasr281@gmail.com'
$ws.Range("G88").Value = 'System, dnasr281@gmail.com'
$ws.Range("G89").Value = 'System, dnasr281@gmail.com'
$ws.Range("G93").Value = 'System, dnasr281@gmail.com'
$ws.Range("G95").Value = 'System, dnasr281@gmail.com'
$ws.Range("G102").Value = 'System, dnasr281@gmail.com'
$ws.Range("G109").Value = 'backup@backdoor.com, System'
$ws.Range("G110").Value = 'backup@backdoor.com, System'
$ws.Range("G111").Value = 'backup@backdoor.com, System'
$ws.Range("G112").Value = 'System, dnasr281@gmail.com'
$ws.Range("G113").Value = 'System, dnasr281@gmail.com'
$ws.Range("G114").Value = 'System, dnasr281@gmail.com'
$ws.Range("G115").Value = 'System, dnasr281@gmail.com'
$ws.Range("G119").Value = 'System, dnasr281@gmail.com'
$ws.Range("G121").Value = 'System, dnasr281@gmail.com'
$ws.Range("G128").Value = 'System, dnasr281@gmail.com'
$ws.Range("G135").Value = 'backup@backdoor.com, System'
$ws.Range("G136").Value = 'backup@backdoor.com, System'
$ws.Range("G137").Value = 'backup@backdoor.com, System'
$ws.Range("G138").Value = 'System, dnasr281@gmail.com'
$ws.Range("G139").Value = 'System, dnasr281@gmail.com'
$ws.Range("G140").Value = 'System, dnasr281@gmail.com'
$ws.Range("G141").Value = 'System, dnasr281@gmail.com'
$ws.Range("G145").Value = 'System, dnasr281@gmail.com'
$ws.Range("G147").Value = 'System, dnasr281@gmail.com'
$ws.Range("G154").Value = 'System, dnasr281@gmail.com'
